# The slide's "Diagram 3" graphic frame (Shape 1) is a SmartArt process
# diagram (ppt/diagrams/data1.xml). The second node in the flow currently
# reads "Rút các bộ ba quan hệ về từ"; update it to read
# "Rút trích các bộ ba" as in the target revision.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)

$sa = $sh.SmartArt
$node = $sa.Nodes.Item(2)
$node.TextFrame2.TextRange.Text = "Rút trích các bộ ba"
